# Update currentAveragePrice / LevePrice / LeveProfit figures across the
# per-job Leve Profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to
# reflect refreshed market board data, as captured by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 684.625
$ws.Range("I2").Value = 684.625
$ws.Range("K2").Value = 684.625
$ws.Range("M2").Value = -571.625

# Row 103
$ws.Range("H103").Value = 560.2143
$ws.Range("I103").Value = 506
$ws.Range("J103").Value = 657.8
$ws.Range("K103").Value = 1518
$ws.Range("L103").Value = 1973.4
$ws.Range("M103").Value = -932
$ws.Range("N103").Value = -3145.4

# Row 107
$ws.Range("H107").Value = 10419270
$ws.Range("I107").Value = 11364090
$ws.Range("J107").Value = 26250
$ws.Range("K107").Value = 11364090
$ws.Range("L107").Value = 26250
$ws.Range("M107").Value = -11362170
$ws.Range("N107").Value = -30090

# Row 112
$ws.Range("H112").Value = 12196257
$ws.Range("I112").Value = 832.8570999999999
$ws.Range("J112").Value = 14707080
$ws.Range("K112").Value = 2498.5713
$ws.Range("L112").Value = 44121240
$ws.Range("M112").Value = -1390.5713
$ws.Range("N112").Value = -44123456

# Row 115
$ws.Range("H115").Value = 709
$ws.Range("I115").Value = 709
$ws.Range("K115").Value = 2127
$ws.Range("M115").Value = -560

# Row 118
$ws.Range("H118").Value = 878.75
$ws.Range("I118").Value = 579.75
$ws.Range("J118").Value = 1177.75
$ws.Range("K118").Value = 1739.25
$ws.Range("L118").Value = 3533.25
$ws.Range("M118").Value = -82.25
$ws.Range("N118").Value = -6847.25

# Row 123
$ws.Range("H123").Value = 28110.285
$ws.Range("J123").Value = 28110.285
$ws.Range("L123").Value = 28110.285
$ws.Range("N123").Value = -37910.285

# Row 126
$ws.Range("H126").Value = 29508.572
$ws.Range("J126").Value = 29508.572
$ws.Range("L126").Value = 29508.572
$ws.Range("N126").Value = -39388.572

# Row 128
$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960

# Row 130
$ws.Range("H130").Value = 42000
$ws.Range("J130").Value = 42000
$ws.Range("L130").Value = 42000
$ws.Range("N130").Value = -52040

# Row 136
$ws.Range("H136").Value = 21374.875
$ws.Range("J136").Value = 21374.875
$ws.Range("L136").Value = 21374.875
$ws.Range("N136").Value = -31574.875

# Row 139
$ws.Range("H139").Value = 21999.5
$ws.Range("J139").Value = 21999.5
$ws.Range("L139").Value = 21999.5
$ws.Range("N139").Value = -32279.5

# Row 140
$ws.Range("H140").Value = 24500
$ws.Range("J140").Value = 24500
$ws.Range("L140").Value = 24500
$ws.Range("N140").Value = -34860

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 7430.2104
$ws.Range("I61").Value = 7430.2104
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7430.2104
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -7218.2104
$ws.Range("N61").Value = ""

# Row 63
$ws.Range("H63").Value = 111114140
$ws.Range("I63").Value = 500001500
$ws.Range("J63").Value = 3471.2856
$ws.Range("K63").Value = 500001500
$ws.Range("L63").Value = 3471.2856
$ws.Range("M63").Value = -500000814
$ws.Range("N63").Value = -4843.2856

# Row 66
$ws.Range("H66").Value = 111114140
$ws.Range("I66").Value = 500001500
$ws.Range("J66").Value = 3471.2856
$ws.Range("K66").Value = 2500007500
$ws.Range("L66").Value = 17356.428
$ws.Range("M66").Value = -2500004068
$ws.Range("N66").Value = -24220.428

# Row 136
$ws.Range("H136").Value = 7430.2104
$ws.Range("I136").Value = 7430.2104
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 22290.6312
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -19740.6312
$ws.Range("N136").Value = ""

# Row 138
$ws.Range("H138").Value = 57500
$ws.Range("J138").Value = 57500
$ws.Range("L138").Value = 57500
$ws.Range("N138").Value = -67780

# Row 139
$ws.Range("H139").Value = 49000
$ws.Range("J139").Value = 49000
$ws.Range("L139").Value = 49000
$ws.Range("N139").Value = -59280

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = -45260

$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 37244.25
$ws.Range("J81").Value = 37244.25
$ws.Range("L81").Value = 37244.25
$ws.Range("N81").Value = -39366.25

# Row 84
$ws.Range("H84").Value = 37244.25
$ws.Range("J84").Value = 37244.25
$ws.Range("L84").Value = 111732.75
$ws.Range("N84").Value = -122340.75

# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""

# Row 137
$ws.Range("H137").Value = 42740
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 42740
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 42740
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -52940

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

# Row 140
$ws.Range("H140").Value = 42900
$ws.Range("J140").Value = 42900
$ws.Range("L140").Value = 42900
$ws.Range("N140").Value = -53260

$ws = $wb.Worksheets.Item("CRP")
# Row 138
$ws.Range("H138").Value = 37799.8
$ws.Range("J138").Value = 37799.8
$ws.Range("L138").Value = 37799.8
$ws.Range("N138").Value = -48079.8

# Row 140
$ws.Range("H140").Value = 24532.072
$ws.Range("J140").Value = 24532.072
$ws.Range("L140").Value = 24532.072
$ws.Range("N140").Value = -34892.072

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 384.8125
$ws.Range("I107").Value = 380
$ws.Range("J107").Value = 385.92307
$ws.Range("K107").Value = 1140
$ws.Range("L107").Value = 1157.76921
$ws.Range("M107").Value = 780
$ws.Range("N107").Value = -4997.76921

$ws = $wb.Worksheets.Item("GSM")
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

# Row 140
$ws.Range("H140").Value = 38368.42
$ws.Range("J140").Value = 38368.42
$ws.Range("L140").Value = 38368.42
$ws.Range("N140").Value = -48728.42

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1020.625
$ws.Range("I16").Value = 610.8333
$ws.Range("J16").Value = 2250
$ws.Range("K16").Value = 610.8333
$ws.Range("L16").Value = 2250
$ws.Range("M16").Value = -440.8333
$ws.Range("N16").Value = -2590

# Row 46
$ws.Range("H46").Value = 18519680
$ws.Range("I46").Value = 30303968
$ws.Range("K46").Value = 30303968
$ws.Range("M46").Value = -30303780

# Row 93
$ws.Range("H93").Value = 13895338
$ws.Range("I93").Value = 9420.087
$ws.Range("J93").Value = 38462732
$ws.Range("K93").Value = 9420.087
$ws.Range("L93").Value = 38462732
$ws.Range("M93").Value = -8172.087
$ws.Range("N93").Value = -38465228

# Row 140
$ws.Range("H140").Value = 60929
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 60929
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 60929
$ws.Range("M140").Value = ""
$ws.Range("N140").Value = -71289

# Row 141
$ws.Range("H141").Value = 49800
$ws.Range("J141").Value = 49800
$ws.Range("L141").Value = 49800
$ws.Range("N141").Value = -60160

$ws = $wb.Worksheets.Item("WVR")
# Row 137
$ws.Range("H137").Value = 33920
$ws.Range("J137").Value = 33920
$ws.Range("L137").Value = 33920
$ws.Range("N137").Value = -44120

# Row 138
$ws.Range("H138").Value = 26666.666
$ws.Range("J138").Value = 26666.666
$ws.Range("L138").Value = 26666.666
$ws.Range("N138").Value = -36946.666

# Row 140
$ws.Range("H140").Value = 34800
$ws.Range("J140").Value = 34800
$ws.Range("L140").Value = 34800
$ws.Range("N140").Value = -45160

# Row 141
$ws.Range("H141").Value = 46535
$ws.Range("I141").Value = 19800
$ws.Range("J141").Value = 49876.875
$ws.Range("K141").Value = 19800
$ws.Range("L141").Value = 49876.875
$ws.Range("M141").Value = -14620
$ws.Range("N141").Value = -60236.875
